$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1, matching the style of the existing header row (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), for rows 2 through 79
$iVals = @(5,4,6,9,8,6,2,5,4,5,3,9,5,7,2,5,8,6,7,4,8,5,4,4,6,6,6,8,7,8,6,8,8,2,9,8,9,7,7,6,5,6,1,8,8,8,8,7,8,6,7,8,8,8,8,8,8,9,7,8,6,6,9,4,7,6,5,8,7,9,6,6,8,4,7,6,5,3)
$jVals = @(6,5,6,9,8,6,2,5,4,6,3,9,5,7,3,6,8,6,7,5,8,6,5,5,7,6,7,8,8,8,7,8,8,3,9,8,9,7,7,7,6,6,1,8,8,8,9,8,8,6,8,8,8,8,8,8,8,9,7,8,6,6,10,4,7,6,5,8,7,9,6,6,8,5,7,6,5,3)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
